$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results - 2 & resume_only+")
$ws.Activate()

# The document-property checker used to report the HYPERLINK property as
# missing (a standalone "Property HYPERLINK does not exist!" row followed by
# a failing "Personal Website / HYPERLINK= / 0" result row). Now that the
# checker also looks up the hyperlink's target, that failure row collapses
# into a single passing result row. Deleting the old "does not exist!" row
# shifts every following row up by one, turning the former row 11 into the
# new row 10.
$ws.Rows.Item(10).Delete()

# Fill in the discovered hyperlink target and flip the pass/fail flag.
$ws.Range("C10").Value = "HYPERLINK=http://www.teachmartin.tk/"
$ws.Range("D10").Value = 1

# Restore the view state captured with the sheet (scroll position / selection).
$win = $excel.ActiveWindow
$win.ScrollRow = 32
$win.ScrollColumn = 1
$ws.Range("C46").Select()
